$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet 1 - Initial")
$ws2 = $wb.Worksheets.Item("Sheet 1 - Product Burndown")

# --- Data edits ---
# Sheet 1 - Initial: Planned Hours (B3) reduced from 5 to 0
$ws1.Range("B3").Value = 0

# Sheet 1 - Product Burndown: accomplished SP / actual hours updates
$ws2.Range("D25").Value = 3
$ws2.Range("E26").Value = 6
$ws2.Range("B27").Value = 7

# --- Selection / view updates ---
# Update the selection remembered on "Sheet 1 - Initial" (A3 -> B3)
$ws1.Activate()
$ws1.Range("B3").Select()

# Restore "Sheet 1 - Product Burndown" as the active sheet and update its selection
$ws2.Activate()
$ws2.Range("B27").Select()
